$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "NA" values under the duplicate_image_filename column (column E)
# for rows 2 through 21.
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}

# Guard against the empty shared-string cell F1 (header row, unused column)
# being reinterpreted during load/save; keep it blank as in the original.
$ws.Cells.Item(1, 6).Value = ""
